$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (date-like text) to Text format so values like
# "01-03-2018" are stored as literal shared-string text instead of
# being auto-parsed into date serial numbers.
$colA = $ws.Range("A87:A100")
$colA.NumberFormat = "@"

# row 87
$ws.Range("A87").Value = "01-03-2018"
$ws.Range("B87").Value = "v1"
$ws.Range("C87").Value = "Day/"
$ws.Range("D87").Value = 1500.0
$ws.Range("E87").Value = 1516.0
$ws.Range("F87").Value = 1700.0
$ws.Range("G87").Value = 16.0
$ws.Range("H87").Value = 66.0
$ws.Range("I87").Value = 200.0
$ws.Range("J87").Value = 27200.0
$ws.Range("K87").Value = "Nothing"
$ws.Range("L87").Value = 0.0
$ws.Range("M87").Value = 14000.0
$ws.Range("N87").Value = "nothing"
# row 88
$ws.Range("A88").Value = "02-03-2018"
$ws.Range("B88").Value = "v1"
$ws.Range("C88").Value = "Day/Night"
$ws.Range("D88").Value = 1516.0
$ws.Range("E88").Value = 1520.0
$ws.Range("F88").Value = 3000.0
$ws.Range("G88").Value = 4.0
$ws.Range("H88").Value = 65.0
$ws.Range("I88").Value = 330.0
$ws.Range("J88").Value = 12000.0
$ws.Range("K88").Value = "Engine oil 	        250"
$ws.Range("L88").Value = 1200.0
$ws.Range("M88").Value = -10650.0
$ws.Range("N88").Value = "a"
# row 89
$ws.Range("A89").Value = "03-03-2018"
$ws.Range("B89").Value = "v1"
$ws.Range("C89").Value = "Day/"
$ws.Range("D89").Value = 1530.0
$ws.Range("E89").Value = 1545.0
$ws.Range("F89").Value = 1700.0
$ws.Range("G89").Value = 15.0
$ws.Range("H89").Value = 64.0
$ws.Range("I89").Value = 200.0
$ws.Range("J89").Value = 25500.0
$ws.Range("K89").Value = "Nothing"
$ws.Range("L89").Value = 0.0
$ws.Range("M89").Value = 12700.0
$ws.Range("N89").Value = "o"
# row 90
$ws.Range("A90").Value = "05-03-2018"
$ws.Range("B90").Value = "v1"
$ws.Range("C90").Value = "Day/Night"
$ws.Range("D90").Value = 1545.0
$ws.Range("E90").Value = 1555.0
$ws.Range("F90").Value = 3000.0
$ws.Range("G90").Value = 10.0
$ws.Range("H90").Value = 68.0
$ws.Range("I90").Value = 200.0
$ws.Range("J90").Value = 30000.0
$ws.Range("K90").Value = "Hydraulic oil	        1000"
$ws.Range("L90").Value = 1500.0
$ws.Range("M90").Value = 14900.0
$ws.Range("N90").Value = "pp"
# row 91
$ws.Range("A91").Value = "01-03-2018"
$ws.Range("B91").Value = "v2"
$ws.Range("C91").Value = "Day/Night"
$ws.Range("D91").Value = 1600.0
$ws.Range("E91").Value = 1616.0
$ws.Range("F91").Value = 2000.0
$ws.Range("G91").Value = 16.0
$ws.Range("H91").Value = 66.0
$ws.Range("I91").Value = 150.0
$ws.Range("J91").Value = 32000.0
$ws.Range("K91").Value = "Nothing"
$ws.Range("L91").Value = 0.0
$ws.Range("M91").Value = 22100.0
$ws.Range("N91").Value = "a"
# row 92
$ws.Range("A92").Value = "02-03-2018"
$ws.Range("B92").Value = "v2"
$ws.Range("C92").Value = "Day/"
$ws.Range("D92").Value = 1616.0
$ws.Range("E92").Value = 1625.0
$ws.Range("F92").Value = 1700.0
$ws.Range("G92").Value = 9.0
$ws.Range("H92").Value = 66.0
$ws.Range("I92").Value = 200.0
$ws.Range("J92").Value = 15300.0
$ws.Range("K92").Value = "Hydraulic strainer  250"
$ws.Range("L92").Value = 2000.0
$ws.Range("M92").Value = 100.0
$ws.Range("N92").Value = "cc"
# row 93
$ws.Range("A93").Value = "03-03-2018"
$ws.Range("B93").Value = "v2"
$ws.Range("C93").Value = "Day/Night"
$ws.Range("D93").Value = 1625.0
$ws.Range("E93").Value = 1640.0
$ws.Range("F93").Value = 3000.0
$ws.Range("G93").Value = 15.0
$ws.Range("H93").Value = 65.0
$ws.Range("I93").Value = 200.0
$ws.Range("J93").Value = 45000.0
$ws.Range("K93").Value = "Nothing"
$ws.Range("L93").Value = 0.0
$ws.Range("M93").Value = 32000.0
$ws.Range("N93").Value = "gy"
# row 94
$ws.Range("A94").Value = "06-03-2018"
$ws.Range("B94").Value = "v2"
$ws.Range("C94").Value = "Day/Night"
$ws.Range("D94").Value = 1640.0
$ws.Range("E94").Value = 1650.0
$ws.Range("F94").Value = 3000.0
$ws.Range("G94").Value = 10.0
$ws.Range("H94").Value = 66.0
$ws.Range("I94").Value = 330.0
$ws.Range("J94").Value = 30000.0
$ws.Range("K94").Value = "Track motor oil       1000"
$ws.Range("L94").Value = 1000.0
$ws.Range("M94").Value = 7220.0
$ws.Range("N94").Value = "pp"
# row 95
$ws.Range("A95").Value = "01-03-2018"
$ws.Range("B95").Value = "v3"
$ws.Range("C95").Value = "Day/"
$ws.Range("D95").Value = 1700.0
$ws.Range("E95").Value = 1716.0
$ws.Range("F95").Value = 1700.0
$ws.Range("G95").Value = 16.0
$ws.Range("H95").Value = 60.0
$ws.Range("I95").Value = 220.0
$ws.Range("J95").Value = 27200.0
$ws.Range("K95").Value = "Engine oil 	        250"
$ws.Range("L95").Value = 1500.0
$ws.Range("M95").Value = 12500.0
$ws.Range("N95").Value = "l"
# row 96
$ws.Range("A96").Value = "03-03-2018"
$ws.Range("B96").Value = "v3"
$ws.Range("C96").Value = "Day/Night"
$ws.Range("D96").Value = 1716.0
$ws.Range("E96").Value = 1730.0
$ws.Range("F96").Value = 3000.0
$ws.Range("G96").Value = 14.0
$ws.Range("H96").Value = 70.0
$ws.Range("I96").Value = 330.0
$ws.Range("J96").Value = 42000.0
$ws.Range("K96").Value = "Nothing"
$ws.Range("L96").Value = 1.0
$ws.Range("M96").Value = 18899.0
$ws.Range("N96").Value = "o"
# row 97
$ws.Range("A97").Value = "15-03-2018"
$ws.Range("B97").Value = "v3"
$ws.Range("C97").Value = "Day/Night"
$ws.Range("D97").Value = 1740.0
$ws.Range("E97").Value = 1750.0
$ws.Range("F97").Value = 3000.0
$ws.Range("G97").Value = 10.0
$ws.Range("H97").Value = 65.0
$ws.Range("I97").Value = 220.0
$ws.Range("J97").Value = 30000.0
$ws.Range("K97").Value = "Hydraulic oil	        1000"
$ws.Range("L97").Value = 2000.0
$ws.Range("M97").Value = 13700.0
$ws.Range("N97").Value = "p"
# row 98
$ws.Range("A98").Value = "07-03-2018"
$ws.Range("B98").Value = "v4"
$ws.Range("C98").Value = "Day/Night"
$ws.Range("D98").Value = 1800.0
$ws.Range("E98").Value = 1820.0
$ws.Range("F98").Value = 3000.0
$ws.Range("G98").Value = 20.0
$ws.Range("H98").Value = 63.0
$ws.Range("I98").Value = 150.0
$ws.Range("J98").Value = 60000.0
$ws.Range("K98").Value = "Swing motor oil       500"
$ws.Range("L98").Value = 2000.0
$ws.Range("M98").Value = 48550.0
$ws.Range("N98").Value = "p"
# row 99
$ws.Range("A99").Value = "13-03-2018"
$ws.Range("B99").Value = "v4"
$ws.Range("C99").Value = "Day/Night"
$ws.Range("D99").Value = 1820.0
$ws.Range("E99").Value = 1830.0
$ws.Range("F99").Value = 3000.0
$ws.Range("G99").Value = 10.0
$ws.Range("H99").Value = 61.0
$ws.Range("I99").Value = 220.0
$ws.Range("J99").Value = 30000.0
$ws.Range("K99").Value = "Nothing"
$ws.Range("L99").Value = 0.0
$ws.Range("M99").Value = 16580.0
$ws.Range("N99").Value = "o"
# row 100
$ws.Range("A100").Value = "05-03-2018"
$ws.Range("B100").Value = "v4"
$ws.Range("C100").Value = "Day/Night"
$ws.Range("D100").Value = 1840.0
$ws.Range("E100").Value = 1845.0
$ws.Range("F100").Value = 3000.0
$ws.Range("G100").Value = 5.0
$ws.Range("H100").Value = 69.0
$ws.Range("I100").Value = 100.0
$ws.Range("J100").Value = 15000.0
$ws.Range("K100").Value = "Air filter                 250"
$ws.Range("L100").Value = 1200.0
$ws.Range("M100").Value = 6900.0
$ws.Range("N100").Value = "no"

# Restore column A to the default (General/Normal) style so the saved
# file does not carry a stray text-format style index on these cells,
# matching the plain <c t="s"> (no explicit s=) cells used elsewhere
# in the sheet for shared-string content.
$colA.Style = "Normal"
